# srish_qbr.docx edit:
#  1) Fix typo "cost-efficnet" -> "cost-efficient" in the Use Case line,
#     landing as three runs (the corrected word isolated in its own run).
#  2) "Moreover, there are " (trailing space) -> "Moreover, there are:"
#     landing as two runs (colon isolated in its own run, trailing space
#     dropped).

$d = $word.ActiveDocument

# --- Edit 1: "cost-efficnet" -> "cost-efficient" -------------------------
$rWord = $d.Content
$found1 = $rWord.Find.Execute("efficnet", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "", 0)
if ($found1) {
    $wordStart = $rWord.Start
    $wordEnd = $rWord.End

    # Correct the misspelling in place (keeps it as its own run boundary).
    $rWord.Text = "efficient"

    # Re-select the corrected word range (length grew by 1 character) and
    # force it onto its own run, distinct from the text before/after it.
    $rFixed = $d.Range($wordStart, $wordStart + 9)
    $rFixed.Bold = $true
    $rFixed.Bold = $false
}

# --- Edit 2: "Moreover, there are " -> "Moreover, there are" + ":" -------
$rPara = $d.Content
$found2 = $rPara.Find.Execute("Moreover, there are ", $true, $false, $false, `
                               $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # Drop the trailing space.
    $rPara.Text = "Moreover, there are"
    $afterText = $rPara.End

    # Append the colon as a new, separate run right after it.
    $rColon = $d.Range($afterText, $afterText)
    $rColon.InsertAfter(":")

    $rNewColon = $d.Range($afterText, $afterText + 1)
    $rNewColon.Bold = $true
    $rNewColon.Bold = $false
}
